$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.368.95"
$ws.Range("E2").Value = "  -1.19%  "
$ws.Range("D3").Value = "'2.716.33"
$ws.Range("E3").Value = "  -1.57%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'561.11"
$ws.Range("E5").Value = "  -2.70%  "
$ws.Range("D6").Value = "'156.91"
$ws.Range("E6").Value = "  -1.38%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'0.590"
$ws.Range("E8").Value = "  -2.43%  "
$ws.Range("E9").Value = "  -2.84%  "
$ws.Range("E10").Value = "  +0.38%  "
$ws.Range("D11").Value = "'5.58"
$ws.Range("E11").Value = "  -2.32%  "
$ws.Range("E12").Value = "  -4.35%  "
$ws.Range("D13").Value = "'3.191.59"
$ws.Range("E13").Value = "  -1.77%  "
$ws.Range("D14").Value = "'26.40"
$ws.Range("E14").Value = "  -1.98%  "
$ws.Range("D15").Value = "'63.157.88"
$ws.Range("E15").Value = "  -0.91%  "
$ws.Range("E16").Value = "  -3.22%  "
$ws.Range("D17").Value = "'2.713.83"
$ws.Range("E17").Value = "  -1.79%  "
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").Value = "'4.67"
$ws.Range("E19").Value = "  -4.23%  "
$ws.Range("D20").Value = "'350.99"
$ws.Range("E20").Value = "  -1.95%  "
$ws.Range("D21").Value = "'6.46"
$ws.Range("E21").Value = "  -4.38%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "'0.512"
$ws.Range("E23").Value = "  -4.21%  "
$ws.Range("D24").Value = "'64.11"
$ws.Range("E24").Value = "  -2.20%  "
$ws.Range("E25").Value = "  -1.28%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  -4.85%  "
$ws.Range("D28").Value = "0.0₃0890"
$ws.Range("E28").Value = "  -2.26%  "
$ws.Range("D29").Value = "'1.37"
$ws.Range("E29").Value = "  +9.49%  "
$ws.Range("E30").Value = "  -0.67%  "
$ws.Range("D31").Value = "'7.17"
$ws.Range("E31").Value = "  -1.77%  "
$ws.Range("D32").Value = "'166.02"
$ws.Range("E32").Value = "  -2.03%  "
$ws.Range("E33").Value = "  -0.86%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").Value = "'19.84"
$ws.Range("E35").Value = "  -2.15%  "
$ws.Range("D36").Value = "'4.83"
$ws.Range("E36").Value = "  -2.17%  "
$ws.Range("E37").Value = "  -2.58%  "
$ws.Range("D38").Value = "'345.40"
$ws.Range("E38").Value = "  -0.25%  "
$ws.Range("D39").Value = "'0.965"
$ws.Range("E39").Value = "  -4.10%  "
$ws.Range("D40").Value = "'6.14"
$ws.Range("E40").Value = "  -3.25%  "
$ws.Range("E41").Value = "  -3.99%  "
$ws.Range("D42").Value = "'38.46"
$ws.Range("E42").Value = "  -1.95%  "
$ws.Range("D43").Value = "'21.36"
$ws.Range("E43").Value = "  -2.39%  "
$ws.Range("D44").Value = "'20.67"
$ws.Range("E44").Value = "  -3.67%  "
$ws.Range("D45").Value = "'0.0573"
$ws.Range("E45").Value = "  -3.20%  "
$ws.Range("D46").Value = "'0.624"
$ws.Range("E46").Value = "  -1.39%  "
$ws.Range("D47").Value = "'0.998"
$ws.Range("E47").Value = "  -0.12%  "
$ws.Range("D48").Value = "'131.67"
$ws.Range("E48").Value = "  -2.92%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0246"
$ws.Range("E49").Value = "  -3.67%  "
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").Value = "'11.05"
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("D51").Value = "'0.0983"
$ws.Range("E51").Value = "  -3.77%  "
